# Doing Updates for Financials
# Insert a new "latest year" column (D) into the PES financials sheet,
# shifting the existing D:K data right to E:L, and populate the new
# column D with the new period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at D; existing D:K shift to E:L.
$ws.Columns("D:D").Insert()

# Copy number formats / styles from the (now shifted) old-D column (E)
# into the new D column, row block by row block -- skipping the blank
# separator rows (36 and 78), and rows 5/6 (which never had D:K data),
# so we don't materialize empty cells/rows that shouldn't exist.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Income Statement (rows 5-35) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 590100
$ws.Range("D9").Value = 429900
$ws.Range("D10").Value = 160200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 4400
$ws.Range("D15").Value = 93600
$ws.Range("D17").Value = 599200
$ws.Range("D18").Value = -9100
$ws.Range("D20").Value = 700
$ws.Range("D21").Value = 85200
$ws.Range("D22").Value = 38800
$ws.Range("D23").Value = -47100
$ws.Range("D24").Value = 1900
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -49000
$ws.Range("D27").Value = -49000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -700
$ws.Range("D33").Value = -49000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -49000

# --- Balance Sheet (rows 37-77) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 53600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 130900
$ws.Range("D44").Value = 18900
$ws.Range("D45").Value = 11700
$ws.Range("D46").Value = 215000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 524900
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 741500
$ws.Range("D57").Value = 34100
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 70600
$ws.Range("D60").Value = 104800
$ws.Range("D61").Value = 464600
$ws.Range("D62").Value = 7200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 576500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -388400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 165100
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (rows 79-102) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -49000
$ws.Range("D83").Value = 93600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 39700
$ws.Range("D91").Value = -67100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -60200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -21100

$wb.Save()
